# Auto-generated edit script applying the scheduled-runner profit-sheet refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR 'Profits' worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 105.6
$ws.Range("I38").Value = 105.6
$ws.Range("K38").Value = 316.8
$ws.Range("M38").Value = 55.20000000000005

$ws.Range("H43").Value = 2139.6155
$ws.Range("J43").Value = 1666.6666
$ws.Range("L43").Value = 1666.6666
$ws.Range("N43").Value = -1804.6666

$ws.Range("H76").Value = 3127.875
$ws.Range("I76").Value = 3001
$ws.Range("J76").Value = 3204
$ws.Range("K76").Value = 3001
$ws.Range("L76").Value = 3204
$ws.Range("M76").Value = -2686
$ws.Range("N76").Value = -3834

$ws.Range("H79").Value = 3127.875
$ws.Range("I79").Value = 3001
$ws.Range("J79").Value = 3204
$ws.Range("K79").Value = 3001
$ws.Range("L79").Value = 3204
$ws.Range("M79").Value = -1909
$ws.Range("N79").Value = -5388

$ws.Range("H98").Value = 98872.125
$ws.Range("I98").Value = 1885.7142
$ws.Range("K98").Value = 1885.7142
$ws.Range("M98").Value = -387.7141999999999

$ws.Range("H122").Value = 98872.125
$ws.Range("I122").Value = 1885.7142
$ws.Range("K122").Value = 5657.142599999999
$ws.Range("M122").Value = -3207.142599999999

$ws.Range("H129").Value = 1406.5
$ws.Range("I129").Value = 2532.6
$ws.Range("J129").Value = 1161.6957
$ws.Range("K129").Value = 7597.799999999999
$ws.Range("L129").Value = 3485.0871
$ws.Range("M129").Value = -2597.799999999999
$ws.Range("N129").Value = -13485.0871

$ws.Range("H137").Value = 3594.4902
$ws.Range("I137").Value = 1068.1364
$ws.Range("J137").Value = 5511.0347
$ws.Range("K137").Value = 3204.4092
$ws.Range("L137").Value = 16533.1041
$ws.Range("M137").Value = -654.4092000000001
$ws.Range("N137").Value = -21633.1041

$ws.Range("H138").Value = 1997.0326
$ws.Range("I138").Value = 1045.98
$ws.Range("J138").Value = 3129.238
$ws.Range("K138").Value = 3137.94
$ws.Range("L138").Value = 9387.714
$ws.Range("M138").Value = 2002.06
$ws.Range("N138").Value = -19667.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 70003.336
$ws.Range("J10").Value = 70003.336
$ws.Range("L10").Value = 70003.336
$ws.Range("N10").Value = -70343.336

$ws.Range("H32").Value = 9444.791999999999
$ws.Range("I32").Value = 8380.880999999999
$ws.Range("J32").Value = 23701.2
$ws.Range("K32").Value = 8380.880999999999
$ws.Range("L32").Value = 23701.2
$ws.Range("M32").Value = -8093.880999999999
$ws.Range("N32").Value = -24275.2

$ws.Range("H58").Value = 29990
$ws.Range("J58").Value = 29990
$ws.Range("L58").Value = 29990
$ws.Range("N58").Value = -30850

$ws.Range("H61").Value = 1306.1364
$ws.Range("I61").Value = 1044.6
$ws.Range("J61").Value = 2323.2222
$ws.Range("K61").Value = 1044.6
$ws.Range("L61").Value = 2323.2222
$ws.Range("M61").Value = -832.5999999999999
$ws.Range("N61").Value = -2747.2222

$ws.Range("H74").Value = 1185.6986
$ws.Range("I74").Value = 1080.7333
$ws.Range("J74").Value = 1670.1538
$ws.Range("K74").Value = 1080.7333
$ws.Range("L74").Value = 1670.1538
$ws.Range("M74").Value = -206.7333000000001
$ws.Range("N74").Value = -3418.1538

$ws.Range("H77").Value = 1185.6986
$ws.Range("I77").Value = 1080.7333
$ws.Range("J77").Value = 1670.1538
$ws.Range("K77").Value = 5403.6665
$ws.Range("L77").Value = 8350.769
$ws.Range("M77").Value = -1035.6665
$ws.Range("N77").Value = -17086.769

$ws.Range("H132").Value = 16668809
$ws.Range("I132").Value = 27779184
$ws.Range("J132").Value = 3246.6667
$ws.Range("K132").Value = 83337552
$ws.Range("L132").Value = 9740.000100000001
$ws.Range("M132").Value = -83335022
$ws.Range("N132").Value = -14800.0001

$ws.Range("H136").Value = 1306.1364
$ws.Range("I136").Value = 1044.6
$ws.Range("J136").Value = 2323.2222
$ws.Range("K136").Value = 3133.8
$ws.Range("L136").Value = 6969.6666
$ws.Range("M136").Value = -583.7999999999997
$ws.Range("N136").Value = -12069.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3419.9465
$ws.Range("I134").Value = 2261.111
$ws.Range("J134").Value = 3968.8684
$ws.Range("K134").Value = 6783.333
$ws.Range("L134").Value = 11906.6052
$ws.Range("M134").Value = -4248.333
$ws.Range("N134").Value = -16976.6052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4681.5107
$ws.Range("I31").Value = 2407
$ws.Range("J31").Value = 5014.3657
$ws.Range("K31").Value = 2407
$ws.Range("L31").Value = 5014.3657
$ws.Range("M31").Value = -2112
$ws.Range("N31").Value = -5604.3657

$ws.Range("H34").Value = 4681.5107
$ws.Range("I34").Value = 2407
$ws.Range("J34").Value = 5014.3657
$ws.Range("K34").Value = 2407
$ws.Range("L34").Value = 5014.3657
$ws.Range("M34").Value = -2205
$ws.Range("N34").Value = -5418.3657

$ws.Range("H110").Value = 47999
$ws.Range("J110").Value = 47999
$ws.Range("L110").Value = 47999
$ws.Range("N110").Value = -56179

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2394.5737
$ws.Range("I5").Value = 3347.2424
$ws.Range("J5").Value = 1271.7858
$ws.Range("K5").Value = 10041.7272
$ws.Range("L5").Value = 3815.3574
$ws.Range("M5").Value = -9929.727200000001
$ws.Range("N5").Value = -4039.3574

$ws.Range("H80").Value = 250500370
$ws.Range("I80").Value = 1000390
$ws.Range("J80").Value = 500000350
$ws.Range("K80").Value = 3001170
$ws.Range("L80").Value = 1500001050
$ws.Range("M80").Value = -3000234
$ws.Range("N80").Value = -1500002922

$ws.Range("H83").Value = 250500370
$ws.Range("I83").Value = 1000390
$ws.Range("J83").Value = 500000350
$ws.Range("K83").Value = 9003510
$ws.Range("L83").Value = 4500003150
$ws.Range("M83").Value = -8998830
$ws.Range("N83").Value = -4500012510

$ws.Range("H121").Value = 13737.375
$ws.Range("J121").Value = 18083.166
$ws.Range("L121").Value = 54249.49800000001
$ws.Range("N121").Value = -56869.49800000001

$ws.Range("H123").Value = 10030
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 141086.69
$ws.Range("I125").Value = 752925
$ws.Range("J125").Value = 5122.6113
$ws.Range("K125").Value = 2258775
$ws.Range("L125").Value = 15367.8339
$ws.Range("M125").Value = -2253855
$ws.Range("N125").Value = -25207.8339

$ws.Range("H132").Value = 2709.3865
$ws.Range("I132").Value = 1884.1666
$ws.Range("J132").Value = 3280.6924
$ws.Range("K132").Value = 16957.4994
$ws.Range("L132").Value = 29526.2316
$ws.Range("M132").Value = -14427.4994
$ws.Range("N132").Value = -34586.2316

$ws.Range("H133").Value = 5445
$ws.Range("I133").Value = 5445
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 16335
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -11275
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 2394.5737
$ws.Range("I135").Value = 3347.2424
$ws.Range("J135").Value = 1271.7858
$ws.Range("K135").Value = 30125.1816
$ws.Range("L135").Value = 11446.0722
$ws.Range("M135").Value = -27590.1816
$ws.Range("N135").Value = -16516.0722

$ws.Range("H136").Value = 25003200
$ws.Range("I136").Value = 41668816
$ws.Range("J136").Value = 4775
$ws.Range("K136").Value = 125006448
$ws.Range("L136").Value = 14325
$ws.Range("M136").Value = -125001348
$ws.Range("N136").Value = -24525

$ws.Range("H137").Value = 52642144
$ws.Range("I137").Value = 3407.5
$ws.Range("J137").Value = 142879980
$ws.Range("K137").Value = 10222.5
$ws.Range("L137").Value = 428639940
$ws.Range("M137").Value = -5122.5
$ws.Range("N137").Value = -428650140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 400
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -384

$ws.Range("H80").Value = 3697.4688
$ws.Range("J80").Value = 3434.1177
$ws.Range("L80").Value = 3434.1177
$ws.Range("N80").Value = -5430.1177

$ws.Range("H83").Value = 3697.4688
$ws.Range("J83").Value = 3434.1177
$ws.Range("L83").Value = 17170.5885
$ws.Range("N83").Value = -27154.5885

$ws.Range("H122").Value = 1016.36365
$ws.Range("I122").Value = 760
$ws.Range("J122").Value = 1112.5
$ws.Range("K122").Value = 2280
$ws.Range("L122").Value = 3337.5
$ws.Range("M122").Value = 170
$ws.Range("N122").Value = -8237.5

$ws.Range("H132").Value = 2167.7817
$ws.Range("I132").Value = 1451.0571
$ws.Range("J132").Value = 3422.05
$ws.Range("K132").Value = 4353.1713
$ws.Range("L132").Value = 10266.15
$ws.Range("M132").Value = -1823.1713
$ws.Range("N132").Value = -15326.15

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3529.2083
$ws.Range("I132").Value = 2425.3572
$ws.Range("J132").Value = 5074.6
$ws.Range("K132").Value = 7276.071599999999
$ws.Range("L132").Value = 15223.8
$ws.Range("M132").Value = -4746.071599999999
$ws.Range("N132").Value = -20283.8

$ws.Range("H136").Value = 1169.1765
$ws.Range("I136").Value = 916.5909
$ws.Range("J136").Value = 2756.8572
$ws.Range("K136").Value = 2749.7727
$ws.Range("L136").Value = 8270.571599999999
$ws.Range("M136").Value = -199.7727
$ws.Range("N136").Value = -13370.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5
$ws.Range("K10").Value = 5
$ws.Range("M10").Value = 164

$ws.Range("H109").Value = 39377
$ws.Range("J109").Value = 39377
$ws.Range("L109").Value = 39377
$ws.Range("N109").Value = -42151

$ws.Range("H132").Value = 1498.7778
$ws.Range("I132").Value = 1171.9556
$ws.Range("J132").Value = 3132.889
$ws.Range("K132").Value = 3515.8668
$ws.Range("L132").Value = 9398.667000000001
$ws.Range("M132").Value = -985.8667999999998
$ws.Range("N132").Value = -14458.667

$ws.Range("H136").Value = 230017.73
$ws.Range("I136").Value = 282669.22
$ws.Range("J136").Value = 1861.2222
$ws.Range("K136").Value = 848007.6599999999
$ws.Range("L136").Value = 5583.6666
$ws.Range("M136").Value = -845457.6599999999
$ws.Range("N136").Value = -10683.6666
